$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Remove the now-unneeded trailing row (canton list shrinks from 4 header/label
# rows + 4 data rows down to a single header row + 4 data rows). ---
$ws.Rows.Item(6).Delete()

# --- Row 1: full column header row -----------------------------------------
# A1:E1 use the plain default style (no explicit font/format override).
$ws.Range("E1").ClearFormats()
$ws.Range("G1").ClearFormats()
$ws.Range("I1").ClearFormats()
$ws.Range("J1").ClearFormats()
$ws.Range("K1").ClearFormats()

$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# F1:K1 carry the small (Arial 9pt) unit-label font used elsewhere on the sheet.
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

# --- Row 2: Chancy-Pougny (was row 3) ---------------------------------------
# Font must be applied before NumberFormat so the cell lands on the existing
# shared style (Arial 9pt + "0"/"0.00") instead of spawning a new one.
$ws.Range("A2:E2").Font.Name = "Arial"
$ws.Range("A2:E2").Font.Size = 9

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 509700
$ws.Range("C2").Value = "Chancy-Pougny"
$ws.Range("D2").Value = 1925
$ws.Range("E2").Value = 2013
$ws.Range("F2").Value = 620
$ws.Range("G2").Value = 34.88
$ws.Range("H2").Value = 32.3
$ws.Range("I2").Value = 74.29
$ws.Range("J2").Value = 83.98
$ws.Range("K2").Value = 158.27

$ws.Range("A2").NumberFormat = "0"
$ws.Range("B2").NumberFormat = "0"
$ws.Range("D2").NumberFormat = "0"
$ws.Range("E2").NumberFormat = "0"
$ws.Range("F2:K2").NumberFormat = "0.00"

# --- Row 3: Verbois (was row 4) ---------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 509600
$ws.Range("C3").Value = "Verbois"
$ws.Range("D3").Value = 1943
$ws.Range("E3").Value = 1999
$ws.Range("F3").Value = 620
$ws.Range("G3").Value = 102.8
$ws.Range("H3").Value = 98
$ws.Range("I3").Value = 211
$ws.Range("J3").Value = 255
$ws.Range("K3").Value = 466

# --- Row 4: Seujet (was row 5) ----------------------------------------------
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 509450
$ws.Range("C4").Value = "Seujet"
$ws.Range("D4").Value = 1994
$ws.Range("E4").Clear()
$ws.Range("F4").Value = 405
$ws.Range("G4").Value = 8.7
$ws.Range("H4").Value = 5.6
$ws.Range("I4").Value = 9.8
$ws.Range("J4").Value = 10.2
$ws.Range("K4").Value = 20

# --- Row 5: Chancy-Pougny Dot. passe pois. (was row 6) ----------------------
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 509750
$ws.Range("C5").Value = "Chancy-Pougny Dot. passe pois."
$ws.Range("D5").Value = 2013
$ws.Range("F5").Value = 4.5
$ws.Range("G5").Value = 0.31
$ws.Range("H5").Value = 0.31
$ws.Range("I5").Value = 0.9
$ws.Range("J5").Value = 0.9
$ws.Range("K5").Value = 1.81

$ws.Range("A2:K2").Select()
